$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ngf"
$ws.Range("C2").Value = "Ngfr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.5666645
$ws.Range("H2").Value = 9.133329
$ws.Range("I2").Value = 0.156583237611307
$ws.Range("J2").Value = 0.1473423006975575
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.337665
$ws.Range("N2").Value = 0.67533
$ws.Range("O2").Value = 0.01438349055657064
$ws.Range("P2").Value = 0.0143657952272707
$ws.Range("Q2").Value = 1.5420027683925
$ws.Range("R2").Value = 6.16801107357
$ws.Range("S2").Value = 0.00225221351949949
$ws.Range("T2").Value = 0.002116689320136056

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ngf"
$ws.Range("C3").Value = "Ngfr"
$ws.Range("D3").Value = "Inflammatory-Mac"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.5666645
$ws.Range("H3").Value = 9.133329
$ws.Range("I3").Value = 0.156583237611307
$ws.Range("J3").Value = 0.1473423006975575
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.05783366666666667
$ws.Range("N3").Value = 0.173501
$ws.Range("O3").Value = 0.002463536340313192
$ws.Range("P3").Value = 0.003690758351808291
$ws.Range("Q3").Value = 0.2641069524715
$ws.Range("R3").Value = 1.584641714829
$ws.Range("S3").Value = 0.0003857484961393501
$ws.Range("T3").Value = 0.000543804826874159

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ngf"
$ws.Range("C4").Value = "Ngfr"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.5666645
$ws.Range("H4").Value = 9.133329
$ws.Range("I4").Value = 0.156583237611307
$ws.Range("J4").Value = 0.1473423006975575
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 23.0803745
$ws.Range("N4").Value = 46.160749
$ws.Range("O4").Value = 0.9831529731031161
$ws.Range("P4").Value = 0.981943446420921
$ws.Range("Q4").Value = 105.4003268758552
$ws.Range("R4").Value = 421.601307503421
$ws.Range("S4").Value = 0.1539452755956681
$ws.Range("T4").Value = 0.1446818065505473

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ngf"
$ws.Range("C5").Value = "Ngfr"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.658241666666667
$ws.Range("H5").Value = 10.974725
$ws.Range("I5").Value = 0.1254349480088258
$ws.Range("J5").Value = 0.1770483939670849
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.337665
$ws.Range("N5").Value = 0.67533
$ws.Range("O5").Value = 0.01438349055657064
$ws.Range("P5").Value = 0.0143657952272707
$ws.Range("Q5").Value = 1.235260172375
$ws.Range("R5").Value = 7.411561034249999
$ws.Range("S5").Value = 0.001804192390148874
$ws.Range("T5").Value = 0.002543440973048291

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Ngf"
$ws.Range("C6").Value = "Ngfr"
$ws.Range("D6").Value = "Inflammatory-Mac"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.658241666666667
$ws.Range("H6").Value = 10.974725
$ws.Range("I6").Value = 0.1254349480088258
$ws.Range("J6").Value = 0.1770483939670849
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.05783366666666667
$ws.Range("N6").Value = 0.173501
$ws.Range("O6").Value = 0.002463536340313192
$ws.Range("P6").Value = 0.003690758351808291
$ws.Range("Q6").Value = 0.2115695291361111
$ws.Range("R6").Value = 1.904125762225
$ws.Range("S6").Value = 0.0003090135527650382
$ws.Range("T6").Value = 0.0006534428387082635

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Ngf"
$ws.Range("C7").Value = "Ngfr"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.658241666666667
$ws.Range("H7").Value = 10.974725
$ws.Range("I7").Value = 0.1254349480088258
$ws.Range("J7").Value = 0.1770483939670849
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 23.0803745
$ws.Range("N7").Value = 46.160749
$ws.Range("O7").Value = 0.9831529731031161
$ws.Range("P7").Value = 0.981943446420921
$ws.Range("Q7").Value = 84.43358767817082
$ws.Range("R7").Value = 506.6015260690249
$ws.Range("S7").Value = 0.1233217420659118
$ws.Range("T7").Value = 0.1738515101553284

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Ngf"
$ws.Range("C8").Value = "Ngfr"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 20.939547
$ws.Range("H8").Value = 41.87909399999999
$ws.Range("I8").Value = 0.7179818143798673
$ws.Range("J8").Value = 0.6756093053353576
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.337665
$ws.Range("N8").Value = 0.67533
$ws.Range("O8").Value = 0.01438349055657064
$ws.Range("P8").Value = 0.0143657952272707
$ws.Range("Q8").Value = 7.070552137754999
$ws.Range("R8").Value = 28.28220855101999
$ws.Range("S8").Value = 0.01032708464692227
$ws.Range("T8").Value = 0.009705664934086353

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Ngf"
$ws.Range("C9").Value = "Ngfr"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 20.939547
$ws.Range("H9").Value = 41.87909399999999
$ws.Range("I9").Value = 0.7179818143798673
$ws.Range("J9").Value = 0.6756093053353576
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.05783366666666667
$ws.Range("N9").Value = 0.173501
$ws.Range("O9").Value = 0.002463536340313192
$ws.Range("P9").Value = 0.003690758351808291
$ws.Range("Q9").Value = 1.211010781349
$ws.Range("R9").Value = 7.266064688093999
$ws.Range("S9").Value = 0.001768774291408804
$ws.Range("T9").Value = 0.002493510686225869

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Ngf"
$ws.Range("C10").Value = "Ngfr"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 20.939547
$ws.Range("H10").Value = 41.87909399999999
$ws.Range("I10").Value = 0.7179818143798673
$ws.Range("J10").Value = 0.6756093053353576
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 23.0803745
$ws.Range("N10").Value = 46.160749
$ws.Range("O10").Value = 0.9831529731031161
$ws.Range("P10").Value = 0.981943446420921
$ws.Range("Q10").Value = 483.2925866203514
$ws.Range("R10").Value = 1933.170346481406
$ws.Range("S10").Value = 0.7058859554415362
$ws.Range("T10").Value = 0.6634101297150453


# Remove the now-stale rows 11-13 (sending cluster "Resolving-Mac" removed)
$ws.Range("A11:T13").Delete()
